# chore: update Sheets via scheduled runner
# Refresh computed market/profit figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 53364.42
$ws.Range("I33").Value = 100611.1
$ws.Range("J33").Value = 868.1111
$ws.Range("K33").Value = 100611.1
$ws.Range("L33").Value = 868.1111
$ws.Range("M33").Value = -100382.1
$ws.Range("N33").Value = -1326.1111

$ws.Range("H69").Value = 33621.25
$ws.Range("I69").Value = 29749
$ws.Range("J69").Value = 37493.5
$ws.Range("K69").Value = 89247
$ws.Range("L69").Value = 112480.5
$ws.Range("M69").Value = -88373
$ws.Range("N69").Value = -114228.5

$ws.Range("H72").Value = 33621.25
$ws.Range("I72").Value = 29749
$ws.Range("J72").Value = 37493.5
$ws.Range("K72").Value = 267741
$ws.Range("L72").Value = 337441.5
$ws.Range("M72").Value = -263373
$ws.Range("N72").Value = -346177.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1930220.9
$ws.Range("I74").Value = 927099.9
$ws.Range("K74").Value = 927099.9
$ws.Range("M74").Value = -926225.9

$ws.Range("H77").Value = 1930220.9
$ws.Range("I77").Value = 927099.9
$ws.Range("K77").Value = 4635499.5
$ws.Range("M77").Value = -4631131.5

$ws.Range("H97").Value = 1075.5
$ws.Range("I97").Value = 732.9286
$ws.Range("J97").Value = 1874.8334
$ws.Range("K97").Value = 732.9286
$ws.Range("L97").Value = 1874.8334
$ws.Range("M97").Value = -236.9286
$ws.Range("N97").Value = -2866.8334

$ws.Range("H122").Value = 3720.1667
$ws.Range("I122").Value = 3576.75
$ws.Range("K122").Value = 10730.25
$ws.Range("M122").Value = -8280.25

$ws.Range("H133").Value = 72999.39999999999
$ws.Range("J133").Value = 72999.39999999999
$ws.Range("L133").Value = 72999.39999999999
$ws.Range("N133").Value = -78059.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 10619.8
$ws.Range("I36").Value = 2699.6667
$ws.Range("J36").Value = 22500
$ws.Range("K36").Value = 2699.6667
$ws.Range("L36").Value = 22500
$ws.Range("M36").Value = -2165.6667
$ws.Range("N36").Value = -23568

$ws.Range("H105").Value = 2504.7742
$ws.Range("I105").Value = 2109.1738
$ws.Range("J105").Value = 3642.125
$ws.Range("K105").Value = 2109.1738
$ws.Range("L105").Value = 3642.125
$ws.Range("M105").Value = -362.1738
$ws.Range("N105").Value = -7136.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1324.25
$ws.Range("I16").Value = 999.1
$ws.Range("K16").Value = 999.1
$ws.Range("M16").Value = -712.1

$ws.Range("H58").Value = 1461
$ws.Range("I58").Value = 1506.6875
$ws.Range("J58").Value = 1339.1666
$ws.Range("K58").Value = 1506.6875
$ws.Range("L58").Value = 1339.1666
$ws.Range("M58").Value = -1303.6875
$ws.Range("N58").Value = -1745.1666

$ws.Range("H62").Value = 6412.0835
$ws.Range("I62").Value = 3214.8572
$ws.Range("K62").Value = 3214.8572
$ws.Range("M62").Value = -2590.8572

$ws.Range("H65").Value = 6412.0835
$ws.Range("I65").Value = 3214.8572
$ws.Range("K65").Value = 16074.286
$ws.Range("M65").Value = -12954.286

$ws.Range("H86").Value = 2974.3333
$ws.Range("I86").Value = 2775
$ws.Range("J86").Value = 3074
$ws.Range("K86").Value = 2775
$ws.Range("L86").Value = 3074
$ws.Range("M86").Value = -1652
$ws.Range("N86").Value = -5320

$ws.Range("H89").Value = 2974.3333
$ws.Range("I89").Value = 2775
$ws.Range("J89").Value = 3074
$ws.Range("K89").Value = 13875
$ws.Range("L89").Value = 15370
$ws.Range("M89").Value = -8259
$ws.Range("N89").Value = -26602

$ws.Range("H107").Value = 669.14813
$ws.Range("I107").Value = 649.8333
$ws.Range("K107").Value = 649.8333
$ws.Range("M107").Value = 1270.1667

$ws.Range("H113").Value = 1324.25
$ws.Range("I113").Value = 999.1
$ws.Range("K113").Value = 999.1
$ws.Range("M113").Value = 1170.9

$ws.Range("H132").Value = 9506.448
$ws.Range("I132").Value = 9798.416999999999
$ws.Range("K132").Value = 29395.251
$ws.Range("M132").Value = -26865.251

$ws.Range("H134").Value = 3848106.2
$ws.Range("I134").Value = 1772.3043
$ws.Range("J134").Value = 33336666
$ws.Range("K134").Value = 5316.9129
$ws.Range("L134").Value = 100009998
$ws.Range("M134").Value = -2781.9129
$ws.Range("N134").Value = -100015068

$ws.Range("H136").Value = 1461
$ws.Range("I136").Value = 1506.6875
$ws.Range("J136").Value = 1339.1666
$ws.Range("K136").Value = 4520.0625
$ws.Range("L136").Value = 4017.4998
$ws.Range("M136").Value = -1970.0625
$ws.Range("N136").Value = -9117.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 959.25
$ws.Range("I8").Value = 959.25
$ws.Range("K8").Value = 2877.75
$ws.Range("M8").Value = -2738.75

$ws.Range("H31").Value = 649.5
$ws.Range("I31").Value = 649.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1948.5
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -1660.5

$ws.Range("H114").Value = 3691.25
$ws.Range("J114").Value = 4566.4
$ws.Range("L114").Value = 13699.2
$ws.Range("N114").Value = -20207.2

$ws.Range("H116").Value = 129911.45
$ws.Range("I116").Value = 155225.11
$ws.Range("J116").Value = 16000
$ws.Range("K116").Value = 465675.33
$ws.Range("L116").Value = 48000
$ws.Range("M116").Value = -462233.33
$ws.Range("N116").Value = -54884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7164.852
$ws.Range("I70").Value = 7252.9443
$ws.Range("K70").Value = 7252.9443
$ws.Range("M70").Value = -6982.9443

$ws.Range("H73").Value = 7164.852
$ws.Range("I73").Value = 7252.9443
$ws.Range("K73").Value = 7252.9443
$ws.Range("M73").Value = -6316.9443

$ws.Range("H122").Value = 5505.75
$ws.Range("I122").Value = 8203.5
$ws.Range("K122").Value = 24610.5
$ws.Range("M122").Value = -22160.5

$ws.Range("H132").Value = 2205.2307
$ws.Range("I132").Value = 2234.8333
$ws.Range("J132").Value = 1850
$ws.Range("K132").Value = 6704.499899999999
$ws.Range("L132").Value = 5550
$ws.Range("M132").Value = -4174.499899999999
$ws.Range("N132").Value = -10610

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2489.0454
$ws.Range("J22").Value = 2632.8823
$ws.Range("L22").Value = 2632.8823
$ws.Range("N22").Value = -3222.8823

$ws.Range("H27").Value = 2489.0454
$ws.Range("J27").Value = 2632.8823
$ws.Range("L27").Value = 2632.8823
$ws.Range("N27").Value = -2846.8823

$ws.Range("H40").Value = 6293
$ws.Range("J40").Value = 6989.5
$ws.Range("L40").Value = 6989.5
$ws.Range("N40").Value = -7261.5

$ws.Range("H43").Value = 23755.75
$ws.Range("J43").Value = 23755.75
$ws.Range("L43").Value = 23755.75
$ws.Range("N43").Value = -24141.75

$ws.Range("H93").Value = 1440.5385
$ws.Range("I93").Value = 1522.25
$ws.Range("J93").Value = 1309.8
$ws.Range("K93").Value = 1522.25
$ws.Range("L93").Value = 1309.8
$ws.Range("M93").Value = -274.25
$ws.Range("N93").Value = -3805.8

$ws.Range("H100").Value = 10773.632
$ws.Range("J100").Value = 11237.4375
$ws.Range("L100").Value = 11237.4375
$ws.Range("N100").Value = -12319.4375

$ws.Range("H122").Value = 3011.1538
$ws.Range("I122").Value = 2920.375
$ws.Range("J122").Value = 3156.4
$ws.Range("K122").Value = 8761.125
$ws.Range("L122").Value = 9469.200000000001
$ws.Range("M122").Value = -6311.125
$ws.Range("N122").Value = -14369.2

$ws.Range("H136").Value = 40002900
$ws.Range("J136").Value = 250002400
$ws.Range("L136").Value = 750007200
$ws.Range("N136").Value = -750012300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 28333
$ws.Range("I2").Value = 20000
$ws.Range("K2").Value = 20000
$ws.Range("M2").Value = -19888

$ws.Range("H40").Value = 22996
$ws.Range("J40").Value = 22996
$ws.Range("L40").Value = 22996
$ws.Range("N40").Value = -23294

$ws.Range("H122").Value = 2053.1538
$ws.Range("J122").Value = 2433.5
$ws.Range("L122").Value = 7300.5
$ws.Range("N122").Value = -12200.5
